$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.224.52'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.54%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.069.01'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.51%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '523.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '135.54'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.73%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.069.64'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.39%  '
$ws.Range("E9").Value = '  +4.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.23'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.55%  '
$ws.Range("E11").Value = '  -3.76%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.397'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.26%  '
$ws.Range("E13").Value = '  +1.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.597.72'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.05'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.07%  '
$ws.Range("E16").Value = '  -4.81%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '57.266.67'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.52%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.067.51'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.46%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.28%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.78'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '347.10'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.06%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.23%  '
$ws.Range("E24").Value = '  +0.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.496'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.61%  '
$ws.Range("E26").Value = '  +0.35%  '
$ws.Range("E27").Value = '  -3.10%  '
$ws.Range("E28").Value = '  -10.88%  '
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.09'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -6.56%  '
$ws.Range("E31").Value = '  -3.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.82'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -10.41%  '
$ws.Range("E33").Value = '  -1.09%  '
$ws.Range("E34").Value = '  +0.32%  '
$ws.Range("E35").Value = '  -0.88%  '
$ws.Range("E36").Value = '  -8.16%  '
$ws.Range("E37").Value = '  -5.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.49'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.08%  '
$ws.Range("E39").Value = '  -5.07%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0656'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.80%  '
$ws.Range("E41").Value = '  -6.95%  '
$ws.Range("E42").Value = '  -0.85%  '
$ws.Range("E43").Value = '  -2.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.409.29'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '36.73'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.16%  '
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.107.55'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.40%  '
$ws.Range("E48").Value = '  -2.81%  '
$ws.Range("E49").Value = '  -2.34%  '
$ws.Range("E50").Value = '  -7.91%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.25'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.42%  '
